$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front; existing columns A:X (and their widths,
# header strings and data values) shift right to B:Y.
$ws.Columns("A:A").Insert()

# New header and lab-number value for the inserted column.
$ws.Range("A1").Value = "Lab. #"
$ws.Range("A2").Value = 7184

# Match the original column width style (narrow, "7" character-width column).
# (6.75 is the input value that this engine's pixel-rounded width model maps
# closest to the intended stored width of 7.7109375.)
$ws.Columns("A:A").ColumnWidth = 6.75

# Highlight the whole data row (A2:Y2) with the light-green fill used to mark
# added laboratory numbers.
$ws.Range("A2:Y2").Interior.Color = 12379352
